$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

$ws.Range("A5").Value = "Склад 4"
$ws.Range("A6").Value = "Склад 5"
$ws.Range("A7").Value = "Склад 6"

$ws.Range("A10").Select()
